# Weekly update: insert a new price record at the top of the data (row 9),
# pushing all existing records (old rows 9-41) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 9 - shifts rows 9:41 down to 10:42 and
# extends the used range to A1:T42 (mirrors the row above for formatting,
# e.g. the date style on column D).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with this week's new record.
$ws.Cells.Item(9, 1).Value = 9
$ws.Cells.Item(9, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(9, 3).Value = "Metropolitana"
$ws.Cells.Item(9, 4).Value = 45243
$ws.Cells.Item(9, 5).Value = 13
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100108
$ws.Cells.Item(9, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(9, 9).Value = 100108007
$ws.Cells.Item(9, 10).Value = "Coco"
$ws.Cells.Item(9, 11).Value = "Sin especificar"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 52
$ws.Cells.Item(9, 14).Value = 22000
$ws.Cells.Item(9, 15).Value = 22000
$ws.Cells.Item(9, 16).Value = 22000
$ws.Cells.Item(9, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(9, 18).Value = "Perú"
$ws.Cells.Item(9, 19).Value = 1100
$ws.Cells.Item(9, 20).Value = 20
